# Update "想去人数" (want-to-go count) values for a handful of events
# across the "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) sheets.
# "本地生活" (sheet3) has no matching rows / is unaffected.

$wb = $excel.ActiveWorkbook

$wsExhibit  = $wb.Worksheets.Item("展览")
$wsShow     = $wb.Worksheets.Item("演出")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# "展览" sheet updates (row => new value)
$wsExhibit.Range("F2").Value  = 1340
$wsExhibit.Range("F3").Value  = 1211
$wsExhibit.Range("F6").Value  = 69
$wsExhibit.Range("F7").Value  = 673
$wsExhibit.Range("F11").Value = 2419
$wsExhibit.Range("F15").Value = 244
$wsExhibit.Range("F16").Value = 592
$wsExhibit.Range("F22").Value = 25
$wsExhibit.Range("F24").Value = 4897
$wsExhibit.Range("F25").Value = 215
$wsExhibit.Range("F26").Value = 426
$wsExhibit.Range("F27").Value = 71
$wsExhibit.Range("F39").Value = 1032

# "演出" sheet update
$wsShow.Range("F10").Value = 3

# "全部类型" sheet updates (same events as above, different row numbers)
$wsAllTypes.Range("F2").Value  = 1340
$wsAllTypes.Range("F5").Value  = 1211
$wsAllTypes.Range("F10").Value = 69
$wsAllTypes.Range("F11").Value = 673
$wsAllTypes.Range("F17").Value = 2419
$wsAllTypes.Range("F21").Value = 244
$wsAllTypes.Range("F22").Value = 592
$wsAllTypes.Range("F28").Value = 25
$wsAllTypes.Range("F29").Value = 4897
$wsAllTypes.Range("F30").Value = 215
$wsAllTypes.Range("F31").Value = 426
$wsAllTypes.Range("F32").Value = 71
$wsAllTypes.Range("F42").Value = 1032
$wsAllTypes.Range("F47").Value = 3
